$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the descriptive text: population considered infinite -> finite
$ws.Range("B21").Value = "A população avaliada é considerada finita"

# Update the recalculated statistic values (column B)
$ws.Range("B2").Value  = 217.97142857142856
$ws.Range("B3").Value  = 181.39946248662133
$ws.Range("B4").Value  = 219.20437724555507
$ws.Range("B5").Value  = 9808.714285714286
$ws.Range("B6").Value  = 8958.103703638277
$ws.Range("B7").Value  = 106593.24867790294
$ws.Range("B8").Value  = 4.346980377481731
$ws.Range("B10").Value = 18.902457379466863
$ws.Range("B11").Value = 9.475175228516605
$ws.Range("B12").Value = 86.26819848816466
$ws.Range("B13").Value = 7442.202070393374
$ws.Range("B14").Value = 89.77894561109467
$ws.Range("B15").Value = 18.896238402211203
$ws.Range("B16").Value = 39.57775523772137
$ws.Range("B17").Value = 21.79714285714286
$ws.Range("B19").Value = 0.8444444444444444
$ws.Range("B20").Value = 55.0
$ws.Range("B22").Value = 450.0
